$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cells that already carry the two "value column" styles used
# throughout the sheet: column B (plain wrapped text) and column C (red
# wrapped text). New/edited cells reuse this formatting via copy/paste-
# format so they land on the very same style Excel already has, instead
# of minting a near-duplicate one.
$styleRefB = $ws.Range("B10")
$styleRefC = $ws.Range("C10")

function Apply-ValueStyle($addr, $isColumnC) {
    $r = $ws.Range($addr)
    if ($isColumnC) {
        $styleRefC.Copy()
    } else {
        $styleRefB.Copy()
    }
    $r.PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false
}

function Set-PlainValue($addr, $text, $isColumnC) {
    # Plain text assignment - used for text Excel's COM layer won't try
    # to reinterpret as a number/date.
    $ws.Range($addr).Value = $text
    Apply-ValueStyle $addr $isColumnC
}

function Set-DateLikeText($addr, $text, $isColumnC) {
    # Excel auto-converts strings that look like dates (e.g. "01/01/2023")
    # into date serials when assigned through .Value. Force the cell to
    # Text format first so the literal string is preserved, then restore
    # the normal look (wrap/top-align + font) by pasting the format from
    # an existing reference cell that already carries the right style.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
    Apply-ValueStyle $addr $isColumnC
}

# --- Name: row -----------------------------------------------------------
$ws.Range("B4").Value = "Graduation Monograph II"
$ws.Range("C4").Value = "Graduation Monograph II"

# --- Ativação: row (date-like text -> needs the text-forcing workaround) -
Set-DateLikeText "B8" "01/01/2023" $false
Set-DateLikeText "C8" "01/01/2023" $true

# --- Objetivos: row --------------------------------------------------------
$ws.Range("B10").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C10").Value = "5840730 - Antonio Jefferson da Silva Machado"

# --- Objectives: row (new content in previously empty cells) -------------
Set-PlainValue "B11" "The Graduation Work (TG) aims to integrate, deepen and apply the knowledge acquired throughout the course, preparing and developing the student's ability to perform tasks that are part of the professional performance profile of the physical engineer." $false
Set-PlainValue "C11" "The Graduation Work (TG) aims to integrate, deepen and apply the knowledge acquired throughout the course, preparing and developing the student's ability to perform tasks that are part of the professional performance profile of the physical engineer." $true

# --- Programa resumido: row (date-like replacement text) -----------------
Set-DateLikeText "B13" "01/01/2023" $false
Set-DateLikeText "C13" "01/01/2023" $true

# --- Short syllabus: row (new content) ------------------------------------
Set-PlainValue "B14" "Prepare a monograph of Undergraduate Work under the guidance of a professor and present it to a panel of examiners." $false
Set-PlainValue "C14" "Prepare a monograph of Undergraduate Work under the guidance of a professor and present it to a panel of examiners." $true

# --- Programa: row ----------------------------------------------------------
$ws.Range("B15").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C15").Value = "5840730 - Antonio Jefferson da Silva Machado"

# --- Syllabus: row (new content) --------------------------------------------
Set-PlainValue "B16" "The course program will consist of the following steps: 1) Preparation and writing of a monograph on a previously defined and approved subject in the Undergraduate Work I discipline. 2) Definition and disclosure of the presentation date after delivery of the monograph in advance of at least , 15 working days. 3) Definition of the panel of examiners, consisting of the supervisor and at least two invited professionals, with training in engineering or related areas. 4) Presentation and evaluation of the TG. 5) Publication of the evaluation. In case of approval, the final copy of the monograph (printed and electronic copy) must be delivered with the agreement of the supervisor." $false
Set-PlainValue "C16" "The course program will consist of the following steps: 1) Preparation and writing of a monograph on a previously defined and approved subject in the Undergraduate Work I discipline. 2) Definition and disclosure of the presentation date after delivery of the monograph in advance of at least , 15 working days. 3) Definition of the panel of examiners, consisting of the supervisor and at least two invited professionals, with training in engineering or related areas. 4) Presentation and evaluation of the TG. 5) Publication of the evaluation. In case of approval, the final copy of the monograph (printed and electronic copy) must be delivered with the agreement of the supervisor." $true

# --- Método: row -------------------------------------------------------------
$ws.Range("B18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"

Write-Output "LOM3250 sheet updated"
